$d = $word.ActiveDocument

# --- 1) Insert the new '2016-09-07' date paragraph, right before the final
#        bookmark paragraph (the one holding the _GoBack bookmark). ---
$lastIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($lastIndex)
$bookmarkPara.Range.InsertParagraphBefore()

$dateParaIndex = $lastIndex
$datePara = $d.Paragraphs.Item($dateParaIndex)
$datePara.Range.Text = '2016-09-07 '

# --- 2) Append the day's log entry as a sequence of runs into the bookmark
#        paragraph, right after the existing bookmarkStart/bookmarkEnd. ---
$entryParaIndex = $lastIndex + 1
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('Idag var jag frånvarande i början av dagen på grund av tandläkarbesök, så Hannes och Pontus höll morgonmöte utan mig. De informerade mig sedan om vad jag missat under mötet, sedan fick jag berätta för dem vad jag skulle göra under dagen. Jag började med att göra en ny PHP klass då jag ligger en bit före min planering, jag var ')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('menad att fortsätta med SQL ')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('idag')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter(', enligt min planering men det var redan avklarat. Därför gjorde jag en ')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('Admin')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter(' klass i PHP där jag läste in ')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('datan')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter(' från ett SQL uttryck jag gjort och sedan ”')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('loopade')
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertAfter('” igenom det och skickade vidare det så att Hannes kunde använda det. ')

# --- 3) Two trailing empty paragraphs after the new entry. ---
$p = $d.Paragraphs.Item($entryParaIndex)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($entryParaIndex + 1)
$p.Range.InsertParagraphAfter()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
